# Refresh the cached "last generated on" date stamp that PowerPoint writes
# into the Date placeholders of the slide master, every slide layout, and
# the notes master (Insert > Header & Footer > Apply to All), moving it
# from 25/10/2018 (2018-10-25) to 01/08/2019 (2019-08-01).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -match "^\d\d?[/.]\d\d?[/.]\d{4}$") {
                $tr.Text = $newText
            }
        }
    }
}

# Slide master date placeholder: "10/25/2018" -> "8/1/2019" (en-US, M/D/YYYY)
Update-DatePlaceholder $p.SlideMaster.Shapes "8/1/2019"

# Every slide layout hanging off the master shares the same stamp.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "8/1/2019"
}

# Notes master date placeholder: "25/10/2018" -> "01/08/2019" (en-GB, DD/MM/YYYY)
Update-DatePlaceholder $p.NotesMaster.Shapes "01/08/2019"
